$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.914.97'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.871.09'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.84'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5095'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3668'
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07185'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8911'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07497'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = '1.877.78'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.68'
$ws.Range('E14').Value = '  +5.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.225'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008515'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.16'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9996'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '26.958.92'
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.013'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '2.112.80'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.390'
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.21'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.774'
$ws.Range('E26').Value = '  -3.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.87'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.085'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.59'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.699'
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09158'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05054'
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7486'
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.989'
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('E37').Value = '  +4.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.527'
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5630'
$ws.Range('E39').Value = '  +6.54%  '
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.072'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.627'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.73'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.541'
$ws.Range('E44').Value = '  +3.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1482'
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4778'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9996'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.08'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.557'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.12'
$ws.Range('E51').Value = '  -0.50%  '
